# The CRM batch values recorded on 2021-03-28 included an erroneous
# first reading (old row 24) and the CRM reference values (column C)
# had been entered incorrectly (incrementing instead of the constant
# CRM value). Remove the bad reading and correct the CRM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the erroneous first 2021-03-28 reading; remaining rows shift up
# (old row 25 -> new row 24, ..., old row 30 -> new row 29).
$ws.Rows.Item(24).Delete()

# Correct the CRM (column C) values for the 2021-03-28 / 2021-04-04
# batch rows -- they should all reference the same CRM value.
$ws.Range("C24:C29").Value = 2224.4699999999998

# Update the selection/scroll position to match where the editor left
# off after the edit.
$ws.Range("G28").Select()
